$d = $word.ActiveDocument

# 1. Title (appears twice: main Heading1 title, and bold run near the end)
$d.Content.Find.Execute("Play Lucky Neko: Gigablox for Free - Slot Game Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Lucky Neko: Gigablox for Free", 2)
$d.Content.Find.Execute("Play Lucky Neko: Gigablox for Free - Slot Game Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Lucky Neko: Gigablox for Free", 2)

# 2. "What we like" bullet list items
$d.Content.Find.Execute("Gigablox feature allows for larger symbols and bigger rewards", $true, $false, $false, $false, $false, $true, 1, $false, "Gigablox feature with larger symbols", 2)
$d.Content.Find.Execute("Japanese theme and design with zen-like atmosphere", $true, $false, $false, $false, $false, $true, 1, $false, "Relaxing Japanese theme and design", 2)
$d.Content.Find.Execute("Free Spins feature with blessed symbol and increased paylines", $true, $false, $false, $false, $false, $true, 1, $false, "Calming soundtrack", 2)
$d.Content.Find.Execute("Simple and enjoyable gameplay experience with synced reels", $true, $false, $false, $false, $false, $true, 1, $false, "Free Spins feature with increased reels and paylines", 2)

# 3. "What we don't like" bullet list item
$d.Content.Find.Execute("Low-value symbols lack creativity", $true, $false, $false, $false, $false, $true, 1, $false, "No progressive jackpot feature", 2)

# 4. Meta description (italic run at the end)
$d.Content.Find.Execute("Read a review of Lucky Neko: Gigablox slot game. Play this Japanese-themed slot for free with Gigablox feature and Free Spins with blessed symbol for bigger wins.", $true, $false, $false, $false, $false, $true, 1, $false, "Check out our review of Lucky Neko: Gigablox and play this Japanese-themed slot game for free.", 2)
